$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Series" column header
$ws.Range("D1").Value = "Series"

# Add series values for rows 2-7
$ws.Range("D2").Value = "100-112"
$ws.Range("D3").Value = "200-212"
$ws.Range("D4").Value = "300-312"
$ws.Range("D5").Value = "400-412"
$ws.Range("D6").Value = "500-512"
$ws.Range("D7").Value = "600-612"

# Update selection to match diff (F6)
$ws.Range("F6").Select()
